$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Cluster Name', 'Active cases')
    ,@('3398 BlueCross Elly Kay Mordialloc', 31)
    ,@('3564 Waverley Valley Aged Care Glen Waverley', 16)
    ,@('3601 Baptcare Westhaven community', 13)
    ,@('3646 Mornington Bay Care Community Mount Martha', 17)
    ,@('3647 Aurrum Aged Care Reservoir', 12)
    ,@('3653 Fronditha Thalpori St Albans Aged Care', 20)
    ,@('3975 Aurrum Aged Care Brunswick West', 11)
    ,@('4257 BlueCross The Gables Camberwell', 16)
    ,@('4295 Hope Aged Care Sunshine West', 16)
    ,@('4314 Estia Health Ardeer', 19)
    ,@('44095 Myrniong Primary School', 13)
    ,@('44304 Brighton Primary School Brighton', 29)
    ,@('44414 Buangor Primary School Buangor', 24)
    ,@('44473 Longwarry Primary School', 10)
    ,@('44490 Armadale Primary School Armadale', 28)
    ,@('44584 Badger Creek Primary School BadgerCreek', 12)
    ,@('44593 Torquay P-6 College Torquay', 28)
    ,@('44620 Canterbury Primary School Canterbury', 14)
    ,@('44623 Brunswick North Primary School Brunswick West', 27)
    ,@('44745 Briar Hill Primary School Briar Hill', 18)
    ,@('44761 Coburg North Primary School Coburg', 19)
    ,@('44799 Eastwood Primary School Ringwood East', 25)
    ,@('44828 Cheltenham East Primary SchoolCheltenham', 12)
    ,@('44930 St Albans Heights Primary School', 15)
    ,@('44979 Campbellfield Heights Primary SchoolCampbellfield', 12)
    ,@('45013 Gladstone Views Primary School', 20)
    ,@('45087 Milgate Primary School Doncaster East', 14)
    ,@('45147 Maramba Primary School Narre Warren', 10)
    ,@('45226 Glen Waverley South Primary SchoolGlen Waverley', 12)
    ,@('45257 Roxburgh Rise Primary School RoxburghPark', 16)
    ,@('45719 St Joseph''s Primary School Numurkah', 17)
    ,@('4574 Village Glen Aged Care ResidencesMornington', 10)
    ,@('45764 Our Lady Help of Christian''s PrimarySchool Brunswick East', 10)
    ,@('45808 St Joseph''s Catholic Primary SchoolChelsea', 11)
    ,@('45858 St Bernard''s Primary Coburg', 24)
    ,@('45861 St Oliver Plunkett Primary School PascoeVale', 11)
    ,@('45902 Mother of God Primary School Ardeer', 12)
    ,@('45975 St Thomas More Primary School Hadfield', 12)
    ,@('45988 St Macartan''s Parish Primary School', 12)
    ,@('46074 St Justin''s Catholic Primary SchoolWheelers Hill', 11)
    ,@('46078 Corpus Christi Primary School Werribee', 18)
    ,@('46135 Wesley College Junior School St KildaRoad Melbourne', 11)
    ,@('46208 Mount Scopus Memorial College GandelCampus Burwood', 13)
    ,@('46306 King''s College Warrnambool', 12)
    ,@('46327 Victory Christian College Strathdale', 10)
    ,@('51529 Sirius College Primary School Dallas', 14)
    ,@('Berwick Chase Primary School Berwick', 12)
    ,@('Brandon Park Primary School Wheelers Hill', 11)
    ,@('Gladstone Views Primary School Gladstone Park', 11)
    ,@('JBS Australia Brooklyn', 19)
    ,@('Kororoit Creek Primary School Burnside HeightsOct-Dec', 37)
    ,@('Mambourin Enterprises Allara Deer Park', 31)
    ,@('Oakleigh South Primary School Oakleigh South', 13)
    ,@('PGL Camp Rumbug Foster North', 47)
    ,@('Rosebud Primary School Rosebud', 16)
    ,@('St Christophers Primary School Airport West', 12)
    ,@('St Louis de Montfort''s School Aspendale', 14)
    ,@('Thomastown West Primary School Camp Doxa''sMalmsbury', 10)
    ,@('Torquay Hotel Torquay', 18)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}